# Apply KHL stats update: add 2025-11-02 matches and refresh derived stats
# (as_of_utc -> 2025-11-02T17:00:00Z, build_version 32 -> 34)

$wb = $excel.ActiveWorkbook

# ---- Sheet: Matches_SOG - append the 5 new matches played 2025-11-02 ----
$ws = $wb.Worksheets.Item("Matches_SOG")
$ws.Range("A422").Value = "'897722"
$ws.Range("B422").Value = "2025-11-02T10:00:00"
$ws.Range("C422").Value = "Амур"
$ws.Range("D422").Value = "Адмирал"
$ws.Range("E422").Value = 21
$ws.Range("F422").Value = 30
$ws.Range("G422").Value = "khl_text"

$ws.Range("A423").Value = "'897719"
$ws.Range("B423").Value = "2025-11-02T14:00:00"
$ws.Range("C423").Value = "Трактор"
$ws.Range("D423").Value = "Сибирь"
$ws.Range("E423").Value = 39
$ws.Range("F423").Value = 33
$ws.Range("G423").Value = "khl_text"

$ws.Range("A424").Value = "'897721"
$ws.Range("B424").Value = "2025-11-02T14:30:00"
$ws.Range("C424").Value = "Салават Юлаев"
$ws.Range("D424").Value = "Локомотив"
$ws.Range("E424").Value = 16
$ws.Range("F424").Value = 37
$ws.Range("G424").Value = "khl_text"

$ws.Range("A425").Value = "'897720"
$ws.Range("B425").Value = "2025-11-02T14:30:00"
$ws.Range("C425").Value = "Автомобилист"
$ws.Range("D425").Value = "Торпедо"
$ws.Range("E425").Value = 24
$ws.Range("F425").Value = 24
$ws.Range("G425").Value = "khl_text"

$ws.Range("A426").Value = "'897723"
$ws.Range("B426").Value = "2025-11-02T17:00:00"
$ws.Range("C426").Value = "ХК Сочи"
$ws.Range("D426").Value = "Спартак"
$ws.Range("E426").Value = 30
$ws.Range("F426").Value = 50
$ws.Range("G426").Value = "khl_text"

# ---- Sheet: Shots_HA - refresh as_of_utc for every team row, plus updated totals ----
$ws = $wb.Worksheets.Item("Shots_HA")
$ws.Range("D2").Value = "2025-11-02T17:00:00Z"
$ws.Range("D3").Value = "2025-11-02T17:00:00Z"
$ws.Range("E3").Value = 18
$ws.Range("F3").Value = 24
$ws.Range("G3").Value = 521
$ws.Range("H3").Value = 558
$ws.Range("I3").Value = 28.9
$ws.Range("J3").Value = 31
$ws.Range("D4").Value = "2025-11-02T17:00:00Z"
$ws.Range("F4").Value = 20
$ws.Range("K4").Value = 638
$ws.Range("L4").Value = 561
$ws.Range("M4").Value = 31.9
$ws.Range("N4").Value = 28.1
$ws.Range("D5").Value = "2025-11-02T17:00:00Z"
$ws.Range("D6").Value = "2025-11-02T17:00:00Z"
$ws.Range("E6").Value = 17
$ws.Range("G6").Value = 511
$ws.Range("H6").Value = 594
$ws.Range("I6").Value = 30.1
$ws.Range("J6").Value = 34.9
$ws.Range("D7").Value = "2025-11-02T17:00:00Z"
$ws.Range("D8").Value = "2025-11-02T17:00:00Z"
$ws.Range("D9").Value = "2025-11-02T17:00:00Z"
$ws.Range("D10").Value = "2025-11-02T17:00:00Z"
$ws.Range("D11").Value = "2025-11-02T17:00:00Z"
$ws.Range("D12").Value = "2025-11-02T17:00:00Z"
$ws.Range("F12").Value = 25
$ws.Range("K12").Value = 777
$ws.Range("L12").Value = 608
$ws.Range("M12").Value = 31.1
$ws.Range("N12").Value = 24.3
$ws.Range("D13").Value = "2025-11-02T17:00:00Z"
$ws.Range("D14").Value = "2025-11-02T17:00:00Z"
$ws.Range("D15").Value = "2025-11-02T17:00:00Z"
$ws.Range("D16").Value = "2025-11-02T17:00:00Z"
$ws.Range("E16").Value = 15
$ws.Range("G16").Value = 403
$ws.Range("H16").Value = 418
$ws.Range("I16").Value = 26.9
$ws.Range("J16").Value = 27.9
$ws.Range("D17").Value = "2025-11-02T17:00:00Z"
$ws.Range("D18").Value = "2025-11-02T17:00:00Z"
$ws.Range("F18").Value = 21
$ws.Range("K18").Value = 588
$ws.Range("L18").Value = 650
$ws.Range("M18").Value = 28
$ws.Range("N18").Value = 31
$ws.Range("D19").Value = "2025-11-02T17:00:00Z"
$ws.Range("F19").Value = 13
$ws.Range("K19").Value = 456
$ws.Range("L19").Value = 478
$ws.Range("M19").Value = 35.1
$ws.Range("N19").Value = 36.8
$ws.Range("D20").Value = "2025-11-02T17:00:00Z"
$ws.Range("F20").Value = 25
$ws.Range("K20").Value = 872
$ws.Range("L20").Value = 812
$ws.Range("M20").Value = 34.9
$ws.Range("N20").Value = 32.5
$ws.Range("D21").Value = "2025-11-02T17:00:00Z"
$ws.Range("E21").Value = 17
$ws.Range("G21").Value = 558
$ws.Range("H21").Value = 522
$ws.Range("I21").Value = 32.8
$ws.Range("J21").Value = 30.7
$ws.Range("D22").Value = "2025-11-02T17:00:00Z"
$ws.Range("E22").Value = 17
$ws.Range("G22").Value = 512
$ws.Range("H22").Value = 552
$ws.Range("J22").Value = 32.5
$ws.Range("D23").Value = "2025-11-02T17:00:00Z"

# ---- Sheet: Shots_Summary - refresh as_of_utc for every team row, plus updated totals ----
$ws = $wb.Worksheets.Item("Shots_Summary")
$ws.Range("D2").Value = "2025-11-02T17:00:00Z"
$ws.Range("D3").Value = "2025-11-02T17:00:00Z"
$ws.Range("E3").Value = 42
$ws.Range("F3").Value = 1189
$ws.Range("G3").Value = 1283
$ws.Range("H3").Value = 28.3
$ws.Range("I3").Value = 30.5
$ws.Range("D4").Value = "2025-11-02T17:00:00Z"
$ws.Range("E4").Value = 34
$ws.Range("F4").Value = 1182
$ws.Range("G4").Value = 953
$ws.Range("H4").Value = 34.8
$ws.Range("I4").Value = 28
$ws.Range("D5").Value = "2025-11-02T17:00:00Z"
$ws.Range("D6").Value = "2025-11-02T17:00:00Z"
$ws.Range("E6").Value = 38
$ws.Range("F6").Value = 1103
$ws.Range("G6").Value = 1363
$ws.Range("H6").Value = 29
$ws.Range("I6").Value = 35.9
$ws.Range("D7").Value = "2025-11-02T17:00:00Z"
$ws.Range("D8").Value = "2025-11-02T17:00:00Z"
$ws.Range("D9").Value = "2025-11-02T17:00:00Z"
$ws.Range("D10").Value = "2025-11-02T17:00:00Z"
$ws.Range("D11").Value = "2025-11-02T17:00:00Z"
$ws.Range("D12").Value = "2025-11-02T17:00:00Z"
$ws.Range("E12").Value = 41
$ws.Range("F12").Value = 1278
$ws.Range("G12").Value = 1045
$ws.Range("H12").Value = 31.2
$ws.Range("I12").Value = 25.5
$ws.Range("D13").Value = "2025-11-02T17:00:00Z"
$ws.Range("D14").Value = "2025-11-02T17:00:00Z"
$ws.Range("D15").Value = "2025-11-02T17:00:00Z"
$ws.Range("D16").Value = "2025-11-02T17:00:00Z"
$ws.Range("E16").Value = 39
$ws.Range("F16").Value = 1061
$ws.Range("G16").Value = 1116
$ws.Range("H16").Value = 27.2
$ws.Range("I16").Value = 28.6
$ws.Range("D17").Value = "2025-11-02T17:00:00Z"
$ws.Range("D18").Value = "2025-11-02T17:00:00Z"
$ws.Range("E18").Value = 39
$ws.Range("F18").Value = 1074
$ws.Range("G18").Value = 1344
$ws.Range("H18").Value = 27.5
$ws.Range("I18").Value = 34.5
$ws.Range("D19").Value = "2025-11-02T17:00:00Z"
$ws.Range("E19").Value = 37
$ws.Range("F19").Value = 1310
$ws.Range("G19").Value = 1142
$ws.Range("H19").Value = 35.4
$ws.Range("D20").Value = "2025-11-02T17:00:00Z"
$ws.Range("E20").Value = 45
$ws.Range("F20").Value = 1515
$ws.Range("G20").Value = 1408
$ws.Range("H20").Value = 33.7
$ws.Range("I20").Value = 31.3
$ws.Range("D21").Value = "2025-11-02T17:00:00Z"
$ws.Range("E21").Value = 41
$ws.Range("F21").Value = 1384
$ws.Range("G21").Value = 1306
$ws.Range("H21").Value = 33.8
$ws.Range("I21").Value = 31.9
$ws.Range("D22").Value = "2025-11-02T17:00:00Z"
$ws.Range("E22").Value = 35
$ws.Range("F22").Value = 964
$ws.Range("G22").Value = 1225
$ws.Range("I22").Value = 35
$ws.Range("D23").Value = "2025-11-02T17:00:00Z"

# ---- Sheet: Meta_ext - bump as_of_utc + build_version ----
$ws = $wb.Worksheets.Item("Meta_ext")
$ws.Range("B2").Value = "2025-11-02T17:00:00Z"
$ws.Range("D2").Value = 34
